$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-04-07 Sunday" "2024-04-08 Monday"

Replace-Text "316×4=1264" "487×2=974"
Replace-Text "190×7=1330" "493×9=4437"
Replace-Text "152×5=760" "918×4=3672"
Replace-Text "129×6=774" "496×5=2480"
Replace-Text "212×2=424" "368×3=1104"

Replace-Text "959×8=7672" "360×8=2880"
Replace-Text "514×9=4626" "903×8=7224"
Replace-Text "261×4=1044" "744×9=6696"
Replace-Text "521×4=2084" "972×6=5832"
Replace-Text "120×3=360" "886×8=7088"

Replace-Text "243×2=486" "793×8=6344"
Replace-Text "711×6=4266" "623×2=1246"
Replace-Text "412×4=1648" "633×6=3798"
Replace-Text "670×4=2680" "969×5=4845"
Replace-Text "200×7=1400" "380×2=760"

Replace-Text "967×3=2901" "439×2=878"
Replace-Text "921×3=2763" "787×2=1574"
Replace-Text "900×5=4500" "900×6=5400"
Replace-Text "550×5=2750" "751×5=3755"
Replace-Text "174×6=1044" "432×4=1728"

Replace-Text "332×7=2324" "992×9=8928"
Replace-Text "673×4=2692" "361×4=1444"
Replace-Text "403×4=1612" "807×3=2421"
Replace-Text "991×6=5946" "484×8=3872"
Replace-Text "663×5=3315" "283×8=2264"
